$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet weather -> nans
$ws.Name = "nans"

# Insert a new column at the front (shifts day/temperature/windspeed/event
# headers + data from A:D to B:E), matching the pandas DataFrame export
# that now carries the row index in column A.
$ws.Columns.Item(1).Insert()

# Give the new index column (A2:A14) the same header/border style as the
# rest of the header row, by copying B1's format (reused style, no new
# style entries) onto it, then fill in the 0-based row index values.
$ws.Range("B1").Copy()
$ws.Range("A2:A14").PasteSpecial(-4122)

$indexValues = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)
for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $indexValues[$i]
}

# Re-write the data rows (now columns B:E) with proper numeric types for
# temperature/windspeed, and clear out the sentinel "missing data" labels
# ("not available" / "not measured" / "no event") as well as the "-1"
# windspeed sentinel, so they read back as blank/NaN cells.
#       day          temp  wind  event
$data = @(
    @("1/1/2017",  32,    6,    "Rain"),
    @("1/4/2017",  $null, 9,    "Sunny"),
    @("1/5/2017",  -1,    $null,"Snow"),
    @("1/6/2017",  $null, 7,    $null),
    @("1/7/2017",  32,    $null,"Rain"),
    @("1/8/2017",  $null, $null,"Sunny"),
    @("1/9/2017",  $null, $null,$null),
    @("1/10/2017", 34,    8,    "Cloudy"),
    @("1/11/2017", -4,    $null,"Snow"),
    @("1/12/2017", 26,    12,   "Sunny"),
    @("1/13/2017", 12,    12,   "Rainy"),
    @("1/11/2017", -1,    12,   "Snow"),
    @("1/14/2017", 40,    $null,"Sunny")
)

# Column B holds day strings like "1/1/2017" that must stay literal text
# (not get auto-parsed into date serials). Mark the range as Text first,
# write the values, then strip the formatting back off so the saved
# cells carry no style index at all - matching plain inline strings.
$dayRange = $ws.Range("B2:B14")
$dayRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    $ws.Cells.Item($row, 2).Value = $rowData[0]

    if ($null -eq $rowData[1]) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $rowData[1]
    }

    if ($null -eq $rowData[2]) {
        $ws.Cells.Item($row, 4).ClearContents()
    } else {
        $ws.Cells.Item($row, 4).Value = $rowData[2]
    }

    if ($null -eq $rowData[3]) {
        $ws.Cells.Item($row, 5).ClearContents()
    } else {
        $ws.Cells.Item($row, 5).Value = $rowData[3]
    }
}

$dayRange.ClearFormats()
